# B6-PowerPoint.pptx edit:
#   The three data tables (on the slides that hold the balance-sheet /
#   accounts practice tables) had their table style switched to a new
#   style id ({B22A76B0-80F3-4167-BFF1-2B187990E333}), replacing the
#   previous default ({91CAF516-FD68-4B39-8A9B-94920E296E39}).
#
# We walk every slide, find any shape that carries a table, and re-apply
# the new table style id to it via Table.ApplyStyle - the PowerPoint OM
# equivalent of picking a different style in the Table Design ribbon.

$p = $ppt.ActivePresentation

$newStyleId = "{B22A76B0-80F3-4167-BFF1-2B187990E333}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)

        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId, $false)
        }
    }
}
